$wb = $excel.ActiveWorkbook

# The "DT" sheet is the 5th sheet in the workbook.
$dt = $wb.Worksheets.Item("DT")

# Add header row: Year, Department, ExamTitle, NumSub
$dt.Range("A1").Value = "Year"
$dt.Range("B1").Value = "Department"
$dt.Range("C1").Value = "ExamTitle"
$dt.Range("D1").Value = "NumSub"

# Column B width adjustment (to fit the "Department" header)
$dt.Columns.Item(2).ColumnWidth = 17.5546875

# Make DT the active/selected sheet and set selection to D1
$dt.Activate() | Out-Null
$dt.Range("D1").Select() | Out-Null
